$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Acre -> Tocantins, value 1.09 -> 1.23
$ws.Range("A2").Value = "Tocantins"
$ws.Range("B2").Value = "Diferença 2025/04 - 2025/04"
$ws.Range("C2").Value = 1.23

# Row 3: Rio Grande do Norte -> Piauí, value 1.02 -> 1.12
$ws.Range("A3").Value = "Piauí"
$ws.Range("B3").Value = "Diferença 2025/04 - 2025/04"
$ws.Range("C3").Value = 1.12

# Row 4: Maranhão -> Acre, value 0.97 -> 1.01
$ws.Range("A4").Value = "Acre"
$ws.Range("B4").Value = "Diferença 2025/04 - 2025/04"
$ws.Range("C4").Value = 1.01

# Row 5: Distrito Federal -> Amazonas, value 0.95 -> 0.96
$ws.Range("A5").Value = "Amazonas"
$ws.Range("B5").Value = "Diferença 2025/04 - 2025/04"
$ws.Range("C5").Value = 0.96

# Row 6: Amazonas -> Pará, value 0.9399999999999999 -> 0.93
$ws.Range("A6").Value = "Pará"
$ws.Range("B6").Value = "Diferença 2025/04 - 2025/04"
$ws.Range("C6").Value = 0.93

# Row 7: Roraima -> Alagoas, value 0.9399999999999999 -> 0.91
$ws.Range("A7").Value = "Alagoas"
$ws.Range("B7").Value = "Diferença 2025/04 - 2025/04"
$ws.Range("C7").Value = 0.91

# Row 8: Sergipe (unchanged name), value 0.75 -> 0.88, rank 20º -> 10º
$ws.Range("B8").Value = "Diferença 2025/04 - 2025/04"
$ws.Range("C8").Value = 0.88
$ws.Range("D8").Value = "10º"

# Row 9: Brasil (unchanged name/value), only period text changes
$ws.Range("B9").Value = "Diferença 2025/04 - 2025/04"

# Row 10: Nordeste (unchanged name), value 0.83 -> 0.87
$ws.Range("B10").Value = "Diferença 2025/04 - 2025/04"
$ws.Range("C10").Value = 0.87
